$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BB12").Value = 8661770
$ws.Range("BC12").Value = 8661770
$ws.Range("BD12").Value = "0636ef5852867394ec3e54a8afa578f4"
$ws.Range("BE12").Value = "bc085f44e0290d6909130ffd900181f4"

$ws.Range("BB13").Value = 10896437
$ws.Range("BC13").Value = 10896437
$ws.Range("BD13").Value = "5c54b0f9afccc785170d9ca06e7afff0"
$ws.Range("BE13").Value = "59dfd4fcfd12e17c2096672335956f7d"

$ws.Range("BB14").Value = 5060871
$ws.Range("BC14").Value = 5060871
$ws.Range("BD14").Value = "1dafd4352d7b1880d6878f5885f8aaba"
$ws.Range("BE14").Value = "9fd1f06db52427c4ce462aa8e4e47f27"

$ws.Range("BB15").Value = 6618705
$ws.Range("BC15").Value = 6618705
$ws.Range("BD15").Value = "87e5e1d5cfdd2c8ba72f9283164403ff"
$ws.Range("BE15").Value = "a523929a72814e255bddfead0d6528d6"

$ws.Range("BB16").Value = 7778364
$ws.Range("BC16").Value = 7778364
$ws.Range("BD16").Value = "a4624e19657b1e221f464de9ae8487ea"
$ws.Range("BE16").Value = "b72840cf821f742963c074629f0293c2"

$ws.Range("BB17").Value = 6342961
$ws.Range("BC17").Value = 6342961
$ws.Range("BD17").Value = "36bd7d08f530bf11815be767908844dc"
$ws.Range("BE17").Value = "fabafcac86e29515c4cfc3203b1fd1d3"

$ws.Range("BB18").Value = 7719837
$ws.Range("BC18").Value = 7719837
$ws.Range("BD18").Value = "1c5ea01efecd43463923f6996ac58b7b"
$ws.Range("BE18").Value = "47c59ac78d38568c7a7bd5d6cecc6776"

$ws.Range("BB19").Value = 7173828
$ws.Range("BC19").Value = 7173828
$ws.Range("BD19").Value = "acdfb9f813a1236ad558b1a44d24ef8d"
$ws.Range("BE19").Value = "af10e84d49499252ecd23701cff9ab03"

$ws.Range("BB20").Value = 6425591
$ws.Range("BC20").Value = 6425591
$ws.Range("BD20").Value = "739f102045646ceafd74994d7ad817e6"
$ws.Range("BE20").Value = "fef8f102eea7eb6d388aaae6caac0552"

$ws.Range("BB21").Value = 8214080
$ws.Range("BC21").Value = 8214080
$ws.Range("BD21").Value = "d8b77993889b8b7b38a9c3422155e114"
$ws.Range("BE21").Value = "865136a19942b98742ea5e5d681b03ee"

$ws.Range("BB22").Value = 7350565
$ws.Range("BC22").Value = 7350565
$ws.Range("BD22").Value = "e11fee38f106403b60096dd7e6ed0365"
$ws.Range("BE22").Value = "16379bb903f026bcadeb7990a8eadbdf"

$ws.Range("BB23").Value = 7973836
$ws.Range("BC23").Value = 7973836
$ws.Range("BD23").Value = "b10069c19d0eeb6810340f1edacd5dd3"
$ws.Range("BE23").Value = "ea2fd05e14905888002dc0ffc9ec15bd"

$ws.Range("BB24").Value = 7253125
$ws.Range("BC24").Value = 7253125
$ws.Range("BD24").Value = "509a35321d2f473d85c5aeebec04ba78"
$ws.Range("BE24").Value = "26847b246848562fb4079851726fc9cc"

$ws.Range("BB25").Value = 7765627
$ws.Range("BC25").Value = 7765627
$ws.Range("BD25").Value = "14dc42f8bb66c71c42e54ca14615ab10"
$ws.Range("BE25").Value = "08ab0caab00e4913f39cc7f72076c10b"

$ws.Range("BB26").Value = 8425739
$ws.Range("BC26").Value = 8425739
$ws.Range("BD26").Value = "6e8b3f3def78cca3e73707e4cbc3dd4f"
$ws.Range("BE26").Value = "095e003effc78d1e7de139dd1c1c6d65"

$ws.Range("BB27").Value = 8938420
$ws.Range("BC27").Value = 8938420
$ws.Range("BD27").Value = "a617a75ce575a430b87737eb0c79ad24"
$ws.Range("BE27").Value = "bf40092766421a768e60162e4b9f2780"

$ws.Range("BB28").Value = 8004303
$ws.Range("BC28").Value = 8004303
$ws.Range("BD28").Value = "5df9023d43755d7da14730b104071acf"
$ws.Range("BE28").Value = "e31d94363d917838001c3b2aa9fecb9a"

$ws.Range("BB29").Value = 8373052
$ws.Range("BC29").Value = 8373052
$ws.Range("BD29").Value = "8ebeacd8918d989e7455532c2313f8e5"
$ws.Range("BE29").Value = "528b9da8d0e122be4a6cf9c5e37a80ff"

$ws.Range("BB30").Value = 8380324
$ws.Range("BC30").Value = 8380324
$ws.Range("BD30").Value = "2b8ecfd1f47277d285debe44f50cc254"
$ws.Range("BE30").Value = "8fac7e60e35c2d278563322df43fe94e"

$ws.Range("BB31").Value = 7785774
$ws.Range("BC31").Value = 7785774
$ws.Range("BD31").Value = "66c1a02b76aa49366fbbb4033d938d1e"
$ws.Range("BE31").Value = "504ea208c80af4b8c080fa7cc0680a03"

$ws.Range("BB32").Value = 7590249
$ws.Range("BC32").Value = 7590249
$ws.Range("BD32").Value = "a99e0415a5c3686cf3f339549d4ca6c2"
$ws.Range("BE32").Value = "d71b12045fa10a5b88ae7b19d2e34692"

$ws.Range("BB34").Value = 7232302
$ws.Range("BC34").Value = 7232302
$ws.Range("BD34").Value = "c039a8b29dfc2a98b3b6abae7f7457d4"
$ws.Range("BE34").Value = "eafe11f52e0c90a5d04e8b0ad8f68459"

$ws.Range("BB35").Value = 8242886
$ws.Range("BC35").Value = 8242886
$ws.Range("BD35").Value = "da02dd872eaa353d4affd00e50c22cf4"
$ws.Range("BE35").Value = "1e6a73716100044b8f7555a1e0d137e1"
